$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Set-ParagraphInnerXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    # Exclude the trailing paragraph-mark character so only the paragraph's
    # *content* (runs) is replaced; the w:p / w:pPr themselves are untouched.
    $inner = $d.Range($r.Start, $r.End - 1)
    $frag = "<w:p $wns>" + $innerXml + "</w:p>"
    $inner.InsertXML($frag)
}

# --- "Funktionale Requirements" (Heading2) -> split + proofErr around "Funktionale"
$xml1 = "<w:proofErr w:type='spellStart'/>" + `
    "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Funktionale</w:t></w:r>" + `
    "<w:proofErr w:type='spellEnd'/>" + `
    "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> Requirements</w:t></w:r>"
Set-ParagraphInnerXml 2 $xml1

# --- "1 Veranstaltungs-Verwaltung" -> split + proofErr around "Veranstaltungs-Verwaltung"
$xml2 = "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>1 </w:t></w:r>" + `
    "<w:proofErr w:type='spellStart'/>" + `
    "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Veranstaltungs-Verwaltung</w:t></w:r>" + `
    "<w:proofErr w:type='spellEnd'/>"
Set-ParagraphInnerXml 3 $xml2

# --- "2 Analyse" -> split + proofErr around "Analyse"
$xml3 = "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>2 </w:t></w:r>" + `
    "<w:proofErr w:type='spellStart'/>" + `
    "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Analyse</w:t></w:r>" + `
    "<w:proofErr w:type='spellEnd'/>"
Set-ParagraphInnerXml 9 $xml3

# --- "3 Tweet-Filterung" -> split + proofErr around "Filterung"
$xml4 = "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>3 Tweet-</w:t></w:r>" + `
    "<w:proofErr w:type='spellStart'/>" + `
    "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Filterung</w:t></w:r>" + `
    "<w:proofErr w:type='spellEnd'/>"
Set-ParagraphInnerXml 14 $xml4

Write-Output "text edits done"

$sec = $d.Sections.First
$footer = $sec.Footers.Item(1)
$footer.PageNumbers.Add(2)

Write-Output "footer added"
